$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13, right-hand "Journal de travail" table (columns H:L) was empty;
# fill it in with the new journal entry ("remplissage de la documentation").

# H13: Cours
$ws.Range("H13").Value = "Ict 431"

# I13: Date - copy the date format/style already used by I12 (and the
# rest of the column) so we reuse the existing style instead of minting
# a brand-new number format, then set the actual date value.
$ws.Range("I12").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I13").Value2 = 43558

# J13: Temps
$ws.Range("J13").Value = "90 min "

# K13: actvité
$ws.Range("K13").Value = "théorie"

# L13: Commentaire - long text, so wrap it within the cell.
$ws.Range("L13").Value = "création d'un repository,instalation de gihubdesktop et prise en main"
$ws.Range("L13").WrapText = $true

# Update the active selection to reflect where editing left off.
$ws.Range("M14").Select() | Out-Null
